# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a brand-new "2022-Q1" worksheet (fund-holding detail, same shape
#    as the existing per-quarter sheets) positioned right before "总计".
# 2. Prepend a "2022-Q1" row (15 funds held, 5.11亿元) to the "总计" summary
#    sheet, shifting the existing history rows down by one.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# === 1. New "2022-Q1" detail sheet ========================================

$new = $wb.Worksheets.Add()
$new.Name = "2022-Q1"
$wb.Worksheets.Item("2022-Q1").Move($wb.Worksheets.Item("总计"))
$ws = $wb.Worksheets.Item("2022-Q1")

# Borrow the bold/centered/bordered look of the header row and index column
# from the structurally-identical "2021-Q4" sheet.
$src = $wb.Worksheets.Item("2021-Q4")
$src.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$src.Range("A2:A16").Copy()
$ws.Range("A2:A16").PasteSpecial(-4122)

# Header row
$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# Data rows (leading "'" forces text storage for numeric-looking values,
# same as the other quarterly sheets: fund codes/percentages are text).
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'512980"
$ws.Cells.Item(2,3).Value = "广发中证传媒ETF"
$ws.Cells.Item(2,4).Value = "'44.11"
$ws.Cells.Item(2,5).Value = "'99.38"
$ws.Cells.Item(2,6).Value = "'5.99"
$ws.Cells.Item(2,7).Value = "'2.6422"
$ws.Cells.Item(2,8).Value = 3

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'159869"
$ws.Cells.Item(3,3).Value = "华夏中证动漫游戏ETF"
$ws.Cells.Item(3,4).Value = "'6.20"
$ws.Cells.Item(3,5).Value = "'98.75"
$ws.Cells.Item(3,6).Value = "'12.31"
$ws.Cells.Item(3,7).Value = "'0.7632"
$ws.Cells.Item(3,8).Value = 2

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'516010"
$ws.Cells.Item(4,3).Value = "国泰中证动漫游戏ETF"
$ws.Cells.Item(4,4).Value = "'4.95"
$ws.Cells.Item(4,5).Value = "'98.91"
$ws.Cells.Item(4,6).Value = "'12.02"
$ws.Cells.Item(4,7).Value = "'0.5950"
$ws.Cells.Item(4,8).Value = 2

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'160629"
$ws.Cells.Item(5,3).Value = "鹏华中证传媒指数（LOF）"
$ws.Cells.Item(5,4).Value = "'7.63"
$ws.Cells.Item(5,5).Value = "'92.90"
$ws.Cells.Item(5,6).Value = "'5.29"
$ws.Cells.Item(5,7).Value = "'0.4036"
$ws.Cells.Item(5,8).Value = 3

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'516770"
$ws.Cells.Item(6,3).Value = "华泰柏瑞中证动漫游戏ETF"
$ws.Cells.Item(6,4).Value = "'1.11"
$ws.Cells.Item(6,5).Value = "'96.56"
$ws.Cells.Item(6,6).Value = "'11.96"
$ws.Cells.Item(6,7).Value = "'0.1328"
$ws.Cells.Item(6,8).Value = 2

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'164818"
$ws.Cells.Item(7,3).Value = "工银瑞信中证传媒指数（LOF）A"
$ws.Cells.Item(7,4).Value = "'1.99"
$ws.Cells.Item(7,5).Value = "'92.70"
$ws.Cells.Item(7,6).Value = "'5.55"
$ws.Cells.Item(7,7).Value = "'0.1104"
$ws.Cells.Item(7,8).Value = 3

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'161030"
$ws.Cells.Item(8,3).Value = "富国中证体育产业指数"
$ws.Cells.Item(8,4).Value = "'2.32"
$ws.Cells.Item(8,5).Value = "'93.75"
$ws.Cells.Item(8,6).Value = "'4.73"
$ws.Cells.Item(8,7).Value = "'0.1097"
$ws.Cells.Item(8,8).Value = 5

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'001628"
$ws.Cells.Item(9,3).Value = "招商体育文化休闲股票"
$ws.Cells.Item(9,4).Value = "'2.95"
$ws.Cells.Item(9,5).Value = "'83.21"
$ws.Cells.Item(9,6).Value = "'3.58"
$ws.Cells.Item(9,7).Value = "'0.1056"
$ws.Cells.Item(9,8).Value = 10

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'159805"
$ws.Cells.Item(10,3).Value = "鹏华中证传媒ETF"
$ws.Cells.Item(10,4).Value = "'1.73"
$ws.Cells.Item(10,5).Value = "'96.29"
$ws.Cells.Item(10,6).Value = "'5.75"
$ws.Cells.Item(10,7).Value = "'0.0995"
$ws.Cells.Item(10,8).Value = 3

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "'350002"
$ws.Cells.Item(11,3).Value = "天治低碳经济灵活配置混合"
$ws.Cells.Item(11,4).Value = "'0.76"
$ws.Cells.Item(11,5).Value = "'65.23"
$ws.Cells.Item(11,6).Value = "'6.52"
$ws.Cells.Item(11,7).Value = "'0.0496"
$ws.Cells.Item(11,8).Value = 4

$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "'159725"
$ws.Cells.Item(12,3).Value = "工银瑞信中证线上消费主题交易型开放式指数证券投资基金"
$ws.Cells.Item(12,4).Value = "'0.75"
$ws.Cells.Item(12,5).Value = "'98.18"
$ws.Cells.Item(12,6).Value = "'4.61"
$ws.Cells.Item(12,7).Value = "'0.0346"
$ws.Cells.Item(12,8).Value = 4

$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "'161036"
$ws.Cells.Item(13,3).Value = "富国中证娱乐主题指数增强（LOF）"
$ws.Cells.Item(13,4).Value = "'0.77"
$ws.Cells.Item(13,5).Value = "'93.32"
$ws.Cells.Item(13,6).Value = "'3.47"
$ws.Cells.Item(13,7).Value = "'0.0267"
$ws.Cells.Item(13,8).Value = 6

$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "'002271"
$ws.Cells.Item(14,3).Value = "招商安弘灵活配置混合"
$ws.Cells.Item(14,4).Value = "'0.50"
$ws.Cells.Item(14,5).Value = "'72.34"
$ws.Cells.Item(14,6).Value = "'2.99"
$ws.Cells.Item(14,7).Value = "'0.0150"
$ws.Cells.Item(14,8).Value = 8

$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "'010677"
$ws.Cells.Item(15,3).Value = "工银瑞信中证传媒指数（LOF）C"
$ws.Cells.Item(15,4).Value = "'0.25"
$ws.Cells.Item(15,5).Value = "'92.70"
$ws.Cells.Item(15,6).Value = "'5.55"
$ws.Cells.Item(15,7).Value = "'0.0139"
$ws.Cells.Item(15,8).Value = 3

$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "'516190"
$ws.Cells.Item(16,3).Value = "华夏中证文娱传媒ETF"
$ws.Cells.Item(16,4).Value = "'0.07"
$ws.Cells.Item(16,5).Value = "'96.81"
$ws.Cells.Item(16,6).Value = "'4.91"
$ws.Cells.Item(16,7).Value = "'0.0034"
$ws.Cells.Item(16,8).Value = 4

# === 2. Update the "总计" summary sheet ====================================

$zj = $wb.Worksheets.Item("总计")

# Insert a new row right under the header for the 2022-Q1 summary entry,
# pushing the existing history rows down by one.
$zj.Rows.Item(2).Insert()
# Row-insert copies the header's bold formatting onto the new row by
# default; strip it so the (non-index) data cells start from plain style.
$zj.Range("B2:D2").ClearFormats()

# Give the new row's index cell (A2) the same bold/center/border look as
# the other index-column cells, by copying the format from the row below.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

# Fix up the sequential index column for the pre-existing rows (now 1..5)
$zj.Cells.Item(3,1).Value = 1
$zj.Cells.Item(4,1).Value = 2
$zj.Cells.Item(5,1).Value = 3
$zj.Cells.Item(6,1).Value = 4
$zj.Cells.Item(7,1).Value = 5

# New summary row: 15 funds held, 5.11亿元 market value for 2022-Q1
$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = "2022-Q1"
$zj.Cells.Item(2,3).Value = 15
$zj.Cells.Item(2,4).Value = 5.11
